# Add a new "average transaction pending time(ms)" column, inserted right
# before the existing "average block time(ms)" column in the header row
# (row 1), and append a new data row (row 10) with a full set of results
# including the new metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift header cells R1:Z1 -> S1:AA1 (process right-to-left so we
#        never clobber a cell before it has been read). ---
$headerCols = @("R","S","T","U","V","W","X","Y","Z")
for ($i = $headerCols.Length - 1; $i -ge 0; $i--) {
    $srcCol = $headerCols[$i]
    $dstCol = $headerCols[$i + 1]
    if ($i -eq $headerCols.Length - 1) {
        $dstCol = "AA"
    }
    $ws.Range($dstCol + "1").Value2 = $ws.Range($srcCol + "1").Value2
}

# --- 2. Write the new header text into the now-vacated R1 cell. ---
$ws.Range("R1").Value2 = "average transaction pending time(ms)"

# --- 3. Append the new data row (row 10). ---
$ws.Range("A10").Value2 = 2
$ws.Range("B10").Value2 = 1
$ws.Range("C10").Value2 = 2
$ws.Range("D10").Value2 = 10
$ws.Range("E10").Value2 = 10
$ws.Range("F10").Value2 = 3
$ws.Range("G10").Value2 = 50
$ws.Range("H10").Value2 = 0
$ws.Range("I10").Value2 = 0
$ws.Range("J10").Value2 = 100
$ws.Range("K10").Value2 = 200
$ws.Range("L10").Value2 = $false
$ws.Range("M10").Value2 = 16384
$ws.Range("N10").Value2 = 40
$ws.Range("O10").Value2 = "<-parameter / result->"
$ws.Range("P10").Value2 = 10
$ws.Range("Q10").Value2 = 0.2623854160308838
$ws.Range("R10").Value2 = 111.4573558330536
$ws.Range("S10").Value2 = 26.23854160308838
$ws.Range("T10").Value2 = 0.2623854160308838
$ws.Range("U10").Value2 = 0.0001642704010009766
$ws.Range("V10").Value2 = 0.090625
$ws.Range("W10").Value2 = 0.0453125
$ws.Range("X10").Value2 = 0
$ws.Range("Y10").Value2 = 94.52276921272278
$ws.Range("Z10").Value2 = 12528
$ws.Range("AA10").Value2 = 2320

# --- 4. Match the number formats used by the corresponding cells above
#        (row 1..9) so row 10 renders the same way. ---
$ws.Range("R10:X10").NumberFormat = "0.000000"
$ws.Range("Y10").NumberFormat = "0.000"
$ws.Range("Z10:AA10").NumberFormat = "0.00"
